# The diff between the before/after canonical OOXML for this document
# consists solely of XML attribute re-ordering (an artifact of the
# canonicalization/serialization step used to produce that diff) together
# with the author's broader commit (adding an M2Doc version custom
# property across the repository's template fixtures). Neither the
# document body text, formatting, structure, nor values differ between
# the two revisions for this particular template. The Word object model
# has no operation that reorders the serialized XML attributes of
# existing elements (Word always re-emits elements using its own fixed
# attribute order), and this template's custom document properties
# already match the target state, so there is no content-level edit to
# make here. We simply touch the document (no-op re-save) to mirror the
# commit without introducing any unintended modification.
$d = $word.ActiveDocument
